$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Temps passé" (time spent) and "Reste à faire" (time remaining)
# values for the "Gestion utilisateur" tasks (rows 2-12)

$ws.Range("E2").Value = 60

$ws.Range("E3").Value = 60

$ws.Range("E4").Value = 90
$ws.Range("F4").Value = 0

$ws.Range("E5").Value = 20
$ws.Range("F5").Value = 30

$ws.Range("E6").Value = 20
$ws.Range("F6").Value = 0

$ws.Range("E7").Value = 40
$ws.Range("F7").Value = 10

$ws.Range("E9").Value = 15
$ws.Range("F9").Value = 10

$ws.Range("E12").Value = 30
$ws.Range("F12").Value = 30

# Total "Temps passé" column (mirrors the existing D58 total formula)
$ws.Range("E58").Formula = "=SUM(E2:E57)"

# New summary row: total time spent, expressed in hours
$ws.Range("A61").Value = "Temps passé"
$ws.Range("B61").Formula = "=E58/60"

# Move the visible selection down to the newly added row
$ws.Range("G61").Select()
